# Fix a bug where duplicate "Bigs" (e.g. "David Zhao" / "David Zhao*")
# referenced the same list of "Littles" preferences in memory instead of
# two separate lists. Re-assign the correct "Littles" preference value
# for each row of the pairing table so that each (possibly duplicate)
# "Big" has its own, independent "Littles" value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Corrected pairing data (Big, Little) -- row 1 is the header.
$pairs = @(
    @("Bigs", "Littles"),
    @("David Zhao", "Capitol Hillary Clinton"),
    @("David Zhao*", "Ernie and Bert Sanders"),
    @("Shirali Nigam", "Parry Hotter"),
    @("Disha Jain", "Her Mine E"),
    @("Robyn Guarriello", "Ben Cars- My Luggage"),
    @("Rohni Awasthi", "Larry Richards"),
    @("Disha Jain*", "Ron Ferretly"),
    @("Robyn Guarriello*", "Larry David")
)

for ($i = 0; $i -lt $pairs.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $pairs[$i][0]
    $ws.Cells.Item($row, 2).Value = $pairs[$i][1]
}
